$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2917716402565462
$ws.Range("C2").Value = 1.655778082260271
$ws.Range("D2").Value = 261.3203778131603
$ws.Range("E2").Value = 10.19245300693656
$ws.Range("G2").Value = 273.4603805426137
